$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.326.97"
$ws.Range("E2").Value = "  -5.29%  "

$ws.Range("D3").Value = "1.671.05"
$ws.Range("E3").Value = "  -3.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5090"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.38%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "

$ws.Range("D12").Value = "1.680.97"
$ws.Range("E12").Value = "  -2.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.555"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.75%  "

$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").Value = "1.899.64"
$ws.Range("E15").Value = "  -3.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008529"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -12.77%  "

$ws.Range("D18").Value = "26.391.42"
$ws.Range("E18").Value = "  -5.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.942"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.41%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.207"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.78%  "

$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.665"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1177"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.96%  "

$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05874"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.258"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.325"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.527"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.513"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.639"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.013"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6009"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01619"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.029"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "

$ws.Range("D41").Value = "1.077.75"
$ws.Range("E41").Value = "  -3.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8681"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value = "1.821.63"
$ws.Range("E45").Value = "  -3.24%  "

$ws.Range("E46").Value = "  +3.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.085"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05192"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.35%  "

